$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.322.90'
$ws.Range('E2').Value = '  +2.64%  '
$ws.Range('D3').Value = '3.692.70'
$ws.Range('E3').Value = '  +7.85%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '582.63'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '177.29'
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').Value = '3.687.22'
$ws.Range('E7').Value = '  +7.92%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.614'
$ws.Range('E8').Value = '  +3.66%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.199'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.55'
$ws.Range('E11').Value = '  +21.08%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.608'
$ws.Range('E12').Value = '  +4.43%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '48.98'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000286'
$ws.Range('E14').Value = '  +2.08%  '
$ws.Range('D15').Value = '4.286.89'
$ws.Range('E15').Value = '  +7.90%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '677.01'
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '8.97'
$ws.Range('E17').Value = '  +4.05%  '
$ws.Range('D18').Value = '3.696.95'
$ws.Range('E18').Value = '  +7.96%  '
$ws.Range('D19').Value = '71.382.05'
$ws.Range('E19').Value = '  +2.69%  '
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.91'
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('E23').Value = '  +5.04%  '
$ws.Range('E24').Value = '  +2.76%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '102.26'
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.82'
$ws.Range('E27').Value = '  +5.61%  '
$ws.Range('E28').Value = '  +7.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '35.05'
$ws.Range('E29').Value = '  +4.91%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.14'
$ws.Range('E30').Value = '  +4.68%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.36'
$ws.Range('E31').Value = '  +5.37%  '
$ws.Range('E32').Value = '  +11.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '589.60'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.18'
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('E35').Value = '  +4.84%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '58.94'
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '3.673.85'
$ws.Range('E38').Value = '  +2.39%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.145'
$ws.Range('E39').Value = '  +4.57%  '
$ws.Range('D40').Value = '0.0₃0765'
$ws.Range('E40').Value = '  +5.18%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '35.25'
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.42'
$ws.Range('E42').Value = '  +5.10%  '
$ws.Range('E43').Value = '  +3.84%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0454'
$ws.Range('E44').Value = '  +9.24%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.348'
$ws.Range('E45').Value = '  +4.91%  '
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('E47').Value = '  +8.15%  '
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '136.30'
$ws.Range('E51').Value = '  +3.28%  '
